# Apply "accounting for substitutions" data update to Figure3.xlsx
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Figure3A: update several NS/S counts across populations
# ---------------------------------------------------------------
$wsA = $wb.Worksheets.Item("Figure3A")

$wsA.Range("C4").Value = 6.86754

$wsA.Range("C5").Value = 3.2446299999999999
$wsA.Range("D5").Value = 32.6248
$wsA.Range("E5").Value = 9.7372200000000007

$wsA.Range("C6").Value = 6.86754

$wsA.Range("C7").Value = 14.2118
$wsA.Range("D7").Value = 15.845700000000001

$wsA.Range("C8").Value = 27.209099999999999
$wsA.Range("D8").Value = 45.511299999999999

$wsA.Range("C9").Value = 23.895199999999999
$wsA.Range("D9").Value = 28.180099999999999
$wsA.Range("F9").Value = 0

$wsA.Range("C10").Value = 11.1753
$wsA.Range("D10").Value = 11.485799999999999
$wsA.Range("E10").Value = 16.357299999999999
$wsA.Range("F10").Value = 22.971699999999998
$wsA.Range("G10").Value = 9.4584299999999999
$wsA.Range("H10").Value = 17.2288

$wsA.Range("C11").Value = 1.2848200000000001

# ---------------------------------------------------------------
# Figure3B: recomputed "S/NS evolutionary rate" vs "NS evolutionary
# rate" values; the leftover duplicate G:H block is removed
# ---------------------------------------------------------------
$wsB = $wb.Worksheets.Item("Figure3B")

$wsB.Range("B3").Value = 10.0063
$wsB.Range("C3").Value = 3.2446299999999999

$wsB.Range("B4").Value = 14.5878
$wsB.Range("C4").Value = 14.2118

$wsB.Range("B5").Value = 31.421299999999999
$wsB.Range("C5").Value = 27.209099999999999

$wsB.Range("B6").Value = 24.881399999999999
$wsB.Range("C6").Value = 23.895199999999999

$wsB.Range("B7").Value = 11.2468
$wsB.Range("C7").Value = 11.1753

$wsB.Range("B8").Value = 17.8796
$wsB.Range("C8").Value = 16.357299999999999

$wsB.Range("B9").Value = 11.246700000000001
$wsB.Range("C9").Value = 9.4584299999999999

$wsB.Range("G4:H10").ClearContents()

# ---------------------------------------------------------------
# Figure3C: recomputed "S/NS evolutionary rate" vs "dN/dS" values;
# the leftover duplicate H:I block is removed
# ---------------------------------------------------------------
$wsC = $wb.Worksheets.Item("Figure3C")

$wsC.Range("B3").Value = 10.0063
$wsC.Range("C3").Value = 0.0994529

$wsC.Range("B4").Value = 14.5878
$wsC.Range("C4").Value = 0.89688800000000002

$wsC.Range("B5").Value = 31.421299999999999
$wsC.Range("C5").Value = 0.597854

$wsC.Range("B6").Value = 24.881399999999999
$wsC.Range("C6").Value = 0.84794700000000001

$wsC.Range("B7").Value = 11.2468
$wsC.Range("C7").Value = 0.97296099999999996

$wsC.Range("B8").Value = 17.8796
$wsC.Range("C8").Value = 0.71206199999999997

$wsC.Range("B9").Value = 11.246700000000001
$wsC.Range("C9").Value = 0.54898999999999998

$wsC.Range("H4:I10").ClearContents()

Write-Output "Figure3 data updated"
